$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Sep 26 21:29:44 EDT 2023"
$ws.Range("B3").Value = "Tue Sep 26 21:29:59 EDT 2023"
$ws.Range("B4").Value = "Tue Sep 26 21:30:14 EDT 2023"
